# Append the two newest GSC export rows (2025-12-31 and 2026-01-01) to the
# "Chart" sheet, matching the trailing data already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the first empty row right after the existing data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row1 = $lastRow + 1
$row2 = $lastRow + 2

function Set-TextCell($cell, $text) {
    # Leading apostrophe forces text storage instead of Excel's automatic
    # date-string parsing; ClearFormats() drops the quote-prefix style that
    # gets attached so the cell stays on the default/general style.
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-TextCell $ws.Cells.Item($row1, 1) "2025-12-31"
$ws.Cells.Item($row1, 2).Value = 0
$ws.Cells.Item($row1, 3).Value = 30

Set-TextCell $ws.Cells.Item($row2, 1) "2026-01-01"
$ws.Cells.Item($row2, 2).Value = 0
$ws.Cells.Item($row2, 3).Value = 29
